# Append a new data row (row 59) to each of the 4 worksheets, matching
# the format/style of the existing rows (row 58), and refresh the used
# range / dimension accordingly.

$wb = $excel.ActiveWorkbook

$rows = @(
    @{ Sheet = "ROW35-FE-LIFTER";  A = 45753.86760603009; B = "0x01,0x90"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"; D = "0x01,0x6e"; E = "0xd"; F = 400; G = 568631262647114000000000.0; H = 366; I = 13 },
    @{ Sheet = "ROW35-MID-LIFTER"; A = 45753.71890783565; B = "0x01,0x90"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"; D = "0x01,0x6e"; E = "0xe"; F = 400; G = 568631262647114000000000.0; H = 366; I = 14 },
    @{ Sheet = "ROW02-FE-LIFTER";  A = 45753.85870184028; B = "0x01,0x90"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"; D = "0x01,0x6e"; E = "0x3"; F = 400; G = 568631262647114000000000.0; H = 366; I = 3 },
    @{ Sheet = "ROW02-MID-LIFTER"; A = 45753.9206578588;  B = "0x01,0x90"; C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"; D = "0x01,0x6e"; E = "0x3"; F = 400; G = 985046333984776000000000.0; H = 366; I = 3 }
)

foreach ($rowData in $rows) {
    $ws = $wb.Worksheets.Item($rowData.Sheet)

    # Last currently-used row (should be 58) -> new row is one below it.
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
    $newRow = $lastRow + 1

    # Match the date/time number format used by the column-A cells above.
    $ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($lastRow, 1).NumberFormat

    $ws.Cells.Item($newRow, 1).Value = $rowData.A
    $ws.Cells.Item($newRow, 2).Value = $rowData.B
    $ws.Cells.Item($newRow, 3).Value = $rowData.C
    $ws.Cells.Item($newRow, 4).Value = $rowData.D
    $ws.Cells.Item($newRow, 5).Value = $rowData.E
    $ws.Cells.Item($newRow, 6).Value = $rowData.F
    $ws.Cells.Item($newRow, 7).Value = $rowData.G
    $ws.Cells.Item($newRow, 8).Value = $rowData.H
    $ws.Cells.Item($newRow, 9).Value = $rowData.I
}
